$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("questions")

# Update difficulty level values (column G) to reflect correct difficulty per question
$ws.Range("G1").Value = 3
$ws.Range("G3").Value = 3
$ws.Range("G5").Value = 2

# Update the active selection to match where the edit was made
$ws.Range("G5").Select()
